$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Git Commit ID string used in column AJ (ScriptLatestRunVersion) for rows 2-80
$ws.Range("AJ2:AJ80").Value = "IndicatorQuantiles.R, Git Commit ID: 2e3ff9a54734c37c56b32bb788c6f054c2509b6b"

# Update the pid values in column AH for rows 2-80
$ws.Range("AH2:AH80").Value = 25080
